{"js": "// The document body contains a single 20-row x 5-column table where every\n// cell holds one arithmetic expression (e.g. \"99-12=\"). The edit replaces\n// each expression with a new one, in table (row, column) order, matching\n// the order of the <w:t> runs in the original document.\n//\n// newValues[row][col] is the replacement text for the cell at that\n// position (0-indexed), built directly from the unified diff (old -> new\n// pairs), preserving document order.\nconst newValues = [\n  [\"80-30=\", \"71+24=\", \"45-5=\", \"43-36=\", \"12+45=\"],\n  [\"11+39=\", \"43-24=\", \"98-76=\", \"86-69=\", \"21+64=\"],\n  [\"20-8=\", \"68-57=\", \"90-14=\", \"92-4=\", \"0+34=\"],\n  [\"12+64=\", \"2+79=\", \"8+55=\", \"4+54=\", \"59+30=\"],\n  [\"19+29=\", \"6+34=\", \"63-43=\", \"67-24=\", \"23+1=\"],\n  [\"73-66=\", \"58-29=\", \"11+77=\", \"48-23=\", \"46+38=\"],\n  [\"57-29=\", \"26+44=\", \"54+21=\", \"51+0=\", \"96-27=\"],\n  [\"46-34=\", \"89-71=\", \"69+11=\", \"8+28=\", \"48+50=\"],\n  [\"86+11=\", \"52-41=\", \"49+25=\", \"40+10=\", \"55+39=\"],\n  [\"33+13=\", \"16+68=\", \"94+3=\", \"49-23=\", \"17-4=\"],\n  [\"89-12=\", \"15+13=\", \"42-33=\", \"45+22=\", \"30-21=\"],\n  [\"41+0=\", \"78-19=\", \"58-6=\", \"82-68=\", \"59+38=\"],\n  [\"67-50=\", \"80+16=\", \"85-18=\", \"95-6=\", \"39-7=\"],\n  [\"78-19=\", \"17+12=\", \"42+32=\", \"65-42=\", \"2+79=\"],\n  [\"6+66=\", \"0+81=\", \"35-15=\", \"59-4=\", \"36+18=\"],\n  [\"86-27=\", \"61-45=\", \"47+18=\", \"49-46=\", \"4+47=\"],\n  [\"54+41=\", \"57+36=\", \"19+25=\", \"49+47=\", \"98-77=\"],\n  [\"49-18=\", \"12-1=\", \"52-33=\", \"8+12=\", \"72-47=\"],\n  [\"55-30=\", \"5+87=\", \"58-23=\", \"57-26=\", \"81-27=\"],\n  [\"65+17=\", \"43-9=\", \"84-13=\", \"26+63=\", \"92-15=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    values[r][c] = newValues[r][c];\n  }\n}\n\n// Writing the whole matrix back updates each cell's text run in place,\n// leaving run/paragraph formatting (font, size, alignment) untouched.\ntable.values = values;\nawait context.sync();\n", "ps1": "# The active document's body contains a single 20-row x 5-column table\n# where every cell holds one arithmetic expression (e.g. \"99-12=\"). This\n# replaces each expression with a new one, in table (row, column) order,\n# matching the order of the <w:t> runs in the original document.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"80-30=\", \"71+24=\", \"45-5=\", \"43-36=\", \"12+45=\"),\n    @(\"11+39=\", \"43-24=\", \"98-76=\", \"86-69=\", \"21+64=\"),\n    @(\"20-8=\", \"68-57=\", \"90-14=\", \"92-4=\", \"0+34=\"),\n    @(\"12+64=\", \"2+79=\", \"8+55=\", \"4+54=\", \"59+30=\"),\n    @(\"19+29=\", \"6+34=\", \"63-43=\", \"67-24=\", \"23+1=\"),\n    @(\"73-66=\", \"58-29=\", \"11+77=\", \"48-23=\", \"46+38=\"),\n    @(\"57-29=\", \"26+44=\", \"54+21=\", \"51+0=\", \"96-27=\"),\n    @(\"46-34=\", \"89-71=\", \"69+11=\", \"8+28=\", \"48+50=\"),\n    @(\"86+11=\", \"52-41=\", \"49+25=\", \"40+10=\", \"55+39=\"),\n    @(\"33+13=\", \"16+68=\", \"94+3=\", \"49-23=\", \"17-4=\"),\n    @(\"89-12=\", \"15+13=\", \"42-33=\", \"45+22=\", \"30-21=\"),\n    @(\"41+0=\", \"78-19=\", \"58-6=\", \"82-68=\", \"59+38=\"),\n    @(\"67-50=\", \"80+16=\", \"85-18=\", \"95-6=\", \"39-7=\"),\n    @(\"78-19=\", \"17+12=\", \"42+32=\", \"65-42=\", \"2+79=\"),\n    @(\"6+66=\", \"0+81=\", \"35-15=\", \"59-4=\", \"36+18=\"),\n    @(\"86-27=\", \"61-45=\", \"47+18=\", \"49-46=\", \"4+47=\"),\n    @(\"54+41=\", \"57+36=\", \"19+25=\", \"49+47=\", \"98-77=\"),\n    @(\"49-18=\", \"12-1=\", \"52-33=\", \"8+12=\", \"72-47=\"),\n    @(\"55-30=\", \"5+87=\", \"58-23=\", \"57-26=\", \"81-27=\"),\n    @(\"65+17=\", \"43-9=\", \"84-13=\", \"26+63=\", \"92-15=\")\n)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Length; $c++) {\n        $cell = $t.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $row[$c]\n    }\n}\n"}
